$d = $word.ActiveDocument

$replacements = @(
    @("76×59=", "53×52="),
    @("52×82=", "46×81="),
    @("16×15=", "44×51="),
    @("59×23=", "38×28="),
    @("21×11=", "53×69="),
    @("46×45=", "91×77="),
    @("78×68=", "22×89="),
    @("20×61=", "29×45="),
    @("31×67=", "46×67="),
    @("11×40=", "42×17="),
    @("69×59=", "52×86="),
    @("19×27=", "34×65="),
    @("84×44=", "30×89="),
    @("73×12=", "90×93="),
    @("99×80=", "19×39="),
    @("48×87=", "24×97="),
    @("48×94=", "52×96="),
    @("93×62=", "97×39="),
    @("12×45=", "53×27="),
    @("78×26=", "32×35="),
    @("13×88=", "20×42="),
    @("59×21=", "57×37="),
    @("42×72=", "16×73="),
    @("40×89=", "93×77="),
    @("31×11=", "13×39=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
